$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dimensionen")

# "MPC3 Plus" Gehäuse-Breite (Gesamt/Breite) wurde von 305 auf 296 mm schmaler ausgelegt.
$ws.Range("F9").Value = 296
